$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.24320417504335
$ws.Range("D2").Value = 11.88889642918189
$ws.Range("E2").Value = 19.19262242179203
$ws.Range("F2").Value = 28.40545309357876
$ws.Range("G2").Value = 26.04254107494291
$ws.Range("H2").Value = 13.62968735419056
$ws.Range("J2").Value = 12.97827851689713
$ws.Range("L2").Value = 10.59000308111557
$ws.Range("M2").Value = 16.71809221142846
$ws.Range("O2").Value = 20.37953145616504
$ws.Range("B3").Value = 17.83824278508632
$ws.Range("D3").Value = 11.9309676415068
$ws.Range("E3").Value = 19.14996887570773
$ws.Range("F3").Value = 28.61078416781668
$ws.Range("G3").Value = 26.15305005312597
$ws.Range("H3").Value = 13.69275343806041
$ws.Range("J3").Value = 12.94386836566156
$ws.Range("L3").Value = 10.24369947362223
$ws.Range("M3").Value = 16.40677172382593
$ws.Range("O3").Value = 20.48416991981414
$ws.Range("B4").Value = 17.58540890654488
$ws.Range("D4").Value = 11.95869770348855
$ws.Range("E4").Value = 19.12687772781009
$ws.Range("F4").Value = 28.74566884030216
$ws.Range("G4").Value = 26.23309376239157
$ws.Range("H4").Value = 13.73435113134262
$ws.Range("J4").Value = 12.92435194663764
$ws.Range("L4").Value = 10.02398921359287
$ws.Range("M4").Value = 16.21239387404703
$ws.Range("O4").Value = 20.55439676050287
$ws.Range("B5").Value = 17.48144561631317
$ws.Range("D5").Value = 11.97047613013346
$ws.Range("E5").Value = 19.118256011363
$ws.Range("F5").Value = 28.80284577887337
$ws.Range("G5").Value = 26.26875406542752
$ws.Range("H5").Value = 13.75202477141058
$ws.Range("J5").Value = 12.91681015458692
$ws.Range("L5").Value = 9.93277575573252
$ws.Range("M5").Value = 16.13245081444735
$ws.Range("O5").Value = 20.58451258227465
$ws.Range("B6").Value = 17.46412998307131
$ws.Range("D6").Value = 11.97246083737337
$ws.Range("E6").Value = 19.11687223513549
$ws.Range("F6").Value = 28.81247331764443
$ws.Range("G6").Value = 26.27485843371378
$ws.Range("H6").Value = 13.75500306839653
$ws.Range("J6").Value = 12.91558285966623
$ws.Range("L6").Value = 9.917531498697628
$ws.Range("M6").Value = 16.11913444144905
$ws.Range("O6").Value = 20.58960359399037
$ws.Range("B7").Value = 17.58401042850288
$ws.Range("D7").Value = 11.95885461406578
$ws.Range("E7").Value = 19.12675824978082
$ws.Range("F7").Value = 28.74643100562201
$ws.Range("G7").Value = 26.23356240568613
$ws.Range("H7").Value = 13.7345865607808
$ws.Range("J7").Value = 12.92424856222221
$ws.Range("L7").Value = 10.02276573863663
$ws.Range("M7").Value = 16.21131859394611
$ws.Range("O7").Value = 20.55479685574445
$ws.Range("B8").Value = 18.10451322303162
$ws.Range("D8").Value = 11.90300915375353
$ws.Range("E8").Value = 19.17727799367234
$ws.Range("F8").Value = 28.47441948481815
$ws.Range("G8").Value = 26.0781006103791
$ws.Range("H8").Value = 13.65083546118981
$ws.Range("J8").Value = 12.96608313761584
$ws.Range("L8").Value = 10.47212382213213
$ws.Range("M8").Value = 16.61146471494086
$ws.Range("O8").Value = 20.41436669013957
$ws.Range("B9").Value = 19.0869782919907
$ws.Range("D9").Value = 11.80851784806463
$ws.Range("E9").Value = 19.30051285756823
$ws.Range("F9").Value = 28.01114802929864
$ws.Range("G9").Value = 25.8709524797254
$ws.Range("H9").Value = 13.50944026299988
$ws.Range("J9").Value = 13.06063077096247
$ws.Range("L9").Value = 11.29323725855716
$ws.Range("M9").Value = 17.36725388453422
$ws.Range("O9").Value = 20.18666293807544
$ws.Range("B10").Value = 19.77920188052691
$ws.Range("D10").Value = 11.74819858578922
$ws.Range("E10").Value = 19.4051951636506
$ws.Range("F10").Value = 27.71387113934218
$ws.Range("G10").Value = 25.77955979127447
$ws.Range("H10").Value = 13.4195220208525
$ws.Range("J10").Value = 13.13733624137008
$ws.Range("L10").Value = 11.85529636341033
$ws.Range("M10").Value = 17.9007339249693
$ws.Range("O10").Value = 20.0487645552482
$ws.Range("B11").Value = 20.08649119613278
$ws.Range("D11").Value = 11.72272319660762
$ws.Range("E11").Value = 19.45575269437038
$ws.Range("F11").Value = 27.58806293241884
$ws.Range("G11").Value = 25.75140991807952
$ws.Range("H11").Value = 13.38165732483242
$ws.Range("J11").Value = 13.17371691347814
$ws.Range("L11").Value = 12.10124450328066
$ws.Range("M11").Value = 18.13788142507324
$ws.Range("O11").Value = 19.9924831193159
$ws.Range("B12").Value = 20.20166938665173
$ws.Range("D12").Value = 11.71335790087348
$ws.Range("E12").Value = 19.47530762121273
$ws.Range("F12").Value = 27.54178459625431
$ws.Range("G12").Value = 25.74269618544597
$ws.Range("H12").Value = 13.36775688834087
$ws.Range("J12").Value = 13.18769940785061
$ws.Range("L12").Value = 12.19292244211014
$ws.Range("M12").Value = 18.2268242464324
$ws.Range("O12").Value = 19.97210382119369
$ws.Range("B13").Value = 20.1769178197577
$ws.Range("D13").Value = 11.7153623690567
$ws.Range("E13").Value = 19.47107809076508
$ws.Range("F13").Value = 27.55169074840055
$ws.Range("G13").Value = 25.74448607091255
$ws.Range("H13").Value = 13.37073109497534
$ws.Range("J13").Value = 13.18467899695823
$ws.Range("L13").Value = 12.17324360927861
$ws.Range("M13").Value = 18.20770803158408
$ws.Range("O13").Value = 19.97645126834678
$ws.Range("B14").Value = 20.09599123735369
$ws.Range("D14").Value = 11.72194706710857
$ws.Range("E14").Value = 19.45735334805489
$ws.Range("F14").Value = 27.58422823692294
$ws.Range("G14").Value = 25.75065395210385
$ws.Range("H14").Value = 13.3805049414561
$ws.Range("J14").Value = 13.17486318066875
$ws.Range("L14").Value = 12.10881641210005
$ws.Range("M14").Value = 18.14521638736383
$ws.Range("O14").Value = 19.99078776502741
$ws.Range("B15").Value = 20.04626436356909
$ws.Range("D15").Value = 11.72601704556704
$ws.Range("E15").Value = 19.44899954942381
$ws.Range("F15").Value = 27.60433604529175
$ws.Range("G15").Value = 25.75468580825869
$ws.Range("H15").Value = 13.38654878892062
$ws.Range("J15").Value = 13.16887728911386
$ws.Range("L15").Value = 12.06916153061912
$ws.Range("M15").Value = 18.10682465477736
$ws.Range("O15").Value = 19.99969099362463
$ws.Range("B16").Value = 19.75895873086361
$ws.Range("D16").Value = 11.7499029208355
$ws.Range("E16").Value = 19.40194915581638
$ws.Range("F16").Value = 27.72228326343786
$ws.Range("G16").Value = 25.78167097402044
$ws.Range("H16").Value = 13.42205781809473
$ws.Range("J16").Value = 13.13498792982228
$ws.Range("L16").Value = 11.8390221081823
$ws.Range("M16").Value = 17.88511884734243
$ws.Range("O16").Value = 20.05257295846217
$ws.Range("B17").Value = 19.58069026434478
$ws.Range("D17").Value = 11.76505865297448
$ws.Range("E17").Value = 19.37382878277134
$ws.Range("F17").Value = 27.79705901296431
$ws.Range("G17").Value = 25.80167570961807
$ws.Range("H17").Value = 13.44462068049503
$ws.Range("J17").Value = 13.11457319871666
$ws.Range("L17").Value = 11.69530247513534
$ws.Range("M17").Value = 17.74764512547213
$ws.Range("O17").Value = 20.08667034135635
$ws.Range("B18").Value = 19.47744402368779
$ws.Range("D18").Value = 11.77396074891136
$ws.Range("E18").Value = 19.35793190796894
$ws.Range("F18").Value = 27.84095468089804
$ws.Range("G18").Value = 25.81444474772294
$ws.Range("H18").Value = 13.45788427836605
$ws.Range("J18").Value = 13.10297163225786
$ws.Range("L18").Value = 11.61172535772976
$ws.Range("M18").Value = 17.66805643781365
$ws.Range("O18").Value = 20.10688894126294
$ws.Range("B19").Value = 19.44236742794824
$ws.Range("D19").Value = 11.77700662906985
$ws.Range("E19").Value = 19.35259748450811
$ws.Range("F19").Value = 27.85596906239526
$ws.Range("G19").Value = 25.81898451138435
$ws.Range("H19").Value = 13.46242420101121
$ws.Range("J19").Value = 13.09906790043942
$ws.Range("L19").Value = 11.58327257186407
$ws.Range("M19").Value = 17.64102221535097
$ws.Range("O19").Value = 20.11383863043108
$ws.Range("B20").Value = 19.59974151281718
$ws.Range("D20").Value = 11.76342616645383
$ws.Range("E20").Value = 19.37679364209404
$ws.Range("F20").Value = 27.78900719146109
$ws.Range("G20").Value = 25.7994153525912
$ws.Range("H20").Value = 13.44218921427844
$ws.Range("J20").Value = 13.11673190253584
$ws.Range("L20").Value = 11.71069664538322
$ws.Range("M20").Value = 17.76233346585681
$ws.Range("O20").Value = 20.08297778404037
$ws.Range("B21").Value = 20.11979421008004
$ws.Range("D21").Value = 11.72000534295046
$ws.Range("E21").Value = 19.46137361480258
$ws.Range("F21").Value = 27.57463415097038
$ws.Range("G21").Value = 25.74878936904935
$ws.Range("H21").Value = 13.37762222767329
$ws.Range("J21").Value = 13.17774079911895
$ws.Range("L21").Value = 12.12778021406835
$ws.Range("M21").Value = 18.16359552137883
$ws.Range("O21").Value = 19.98655141376432
$ws.Range("B22").Value = 20.4527304277468
$ws.Range("D22").Value = 11.69326883106
$ws.Range("E22").Value = 19.51903458457909
$ws.Range("F22").Value = 27.44247420050456
$ws.Range("G22").Value = 25.72705082970124
$ws.Range("H22").Value = 13.3379780471521
$ws.Range("J22").Value = 13.2188101152312
$ws.Range("L22").Value = 12.39185251210248
$ws.Range("M22").Value = 18.4208051757893
$ws.Range("O22").Value = 19.92897360752177
$ws.Range("B23").Value = 20.27570039397403
$ws.Range("D23").Value = 11.70738866203236
$ws.Range("E23").Value = 19.48804601073997
$ws.Range("F23").Value = 27.5122810091512
$ws.Range("G23").Value = 25.73761016583937
$ws.Range("H23").Value = 13.3589028540041
$ws.Range("J23").Value = 13.19678383230357
$ws.Range("L23").Value = 12.2517085190633
$ws.Range("M23").Value = 18.28400859873192
$ws.Range("O23").Value = 19.95920407501575
$ws.Range("B24").Value = 19.59113079520747
$ws.Range("D24").Value = 11.76416362525156
$ws.Range("E24").Value = 19.37545238751287
$ws.Range("F24").Value = 27.7926445980482
$ws.Range("G24").Value = 25.80043331130606
$ws.Range("H24").Value = 13.44328757093878
$ws.Range("J24").Value = 13.11575553099716
$ws.Range("L24").Value = 11.70373989811419
$ws.Range("M24").Value = 17.7556945854192
$ws.Range("O24").Value = 20.08464527219237
$ws.Range("B25").Value = 18.82596362338881
$ws.Range("D25").Value = 11.83247770994064
$ws.Range("E25").Value = 19.26464638415603
$ws.Range("F25").Value = 28.1289329054661
$ws.Range("G25").Value = 25.91639674002795
$ws.Range("H25").Value = 13.54524205349108
$ws.Range("J25").Value = 13.03375258465188
$ws.Range("L25").Value = 11.07804388673443
$ws.Range("M25").Value = 17.16633207579774
$ws.Range("O25").Value = 20.24312420227339
